$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.976.34"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.655.31"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("D5").Value = "'215.03"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("D8").Value = "'0.250"
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").Value = "'20.12"
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("E11").Value = "  +3.49%  "
$ws.Range("D12").Value = "1.889.16"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").Value = "1.650.93"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "26.980.50"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "'236.62"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'7.78"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +3.55%  "
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "'145.23"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'7.12"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'15.87"
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").Value = "1.552.36"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("E35").Value = "  +9.58%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'0.581"
$ws.Range("E37").Value = "  +3.54%  "
$ws.Range("D38").Value = "'0.901"
$ws.Range("E38").Value = "  +8.97%  "
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("D40").Value = "'6.03"
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'66.77"
$ws.Range("E42").Value = "  +8.91%  "
$ws.Range("D43").Value = "'0.985"
$ws.Range("E43").Value = "  +5.98%  "
$ws.Range("D44").Value = "'2.24"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").Value = "1.796.97"
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "'89.99"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").Value = "'0.0991"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "'7.69"
$ws.Range("E51").Value = "  +3.46%  "
